# ToDO.xlsx update — "Update ToDo y StackWeapon"
#
# 1. Fill in the previously-blank StackWeapon tracking cells (C7, B8, C8,
#    B9, C9, B10, C10) with their status letters.
# 2. The two "falta" follow-up notes in column G move up two rows
#    (G7->G5, G8->G6); the old G7/G8 cells become empty.
# 3. A thin box-border is added around the whole ToDo table (A1:D10) and
#    around the small PSP0 legend table (G10:J11).
# 4. View settings: zoom to 110% and move the active selection to G12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move the two note cells in column G up two rows -------------------
$ws.Range("G5").Value = $ws.Range("G7").Value2
$ws.Range("G6").Value = $ws.Range("G8").Value2
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()

# --- 2. Fill in the StackWeapon status cells that were still empty --------
$ws.Range("C7").Value = "p"

$ws.Range("B8").Value = "n"
$ws.Range("C8").Value = "p"

$ws.Range("B9").Value = "p"
$ws.Range("C9").Value = "p"

$ws.Range("B10").Value = "n"
$ws.Range("C10").Value = "p"

# --- 3. Add a thin box border around the two tables ------------------------
# Column A (task names) has no alignment style of its own.
$ws.Range("A1:A10").Borders.LineStyle = 1
# Columns B:D (status letters) keep their existing center alignment.
$ws.Range("B1:D10").Borders.LineStyle = 1
# Small legend table bottom-right (names header + PSP0 row).
$ws.Range("G10:J11").Borders.LineStyle = 1

# --- 4. Update the view: zoom 110%, move selection to G12 ------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 110
[void]$ws.Range("G12").Select()
